$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 102, shifting existing rows 102:155 down to 104:157.
$ws.Range("A102:A103").EntireRow.Insert()

# New row 102 data (Vega Modelo de Temuco - Pepino dulce)
$ws.Cells.Item(102, 1).Value = 10
$ws.Cells.Item(102, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value = "La Araucanía"
$ws.Cells.Item(102, 4).Value = 44466
$ws.Cells.Item(102, 5).Value = 9
$ws.Cells.Item(102, 6).Value = 100112043
$ws.Cells.Item(102, 7).Value = "Pepino dulce"
$ws.Cells.Item(102, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(102, 9).Value = "Primera"
$ws.Cells.Item(102, 10).Value = 80
$ws.Cells.Item(102, 11).Value = 24000
$ws.Cells.Item(102, 12).Value = 24000
$ws.Cells.Item(102, 13).Value = 24000
$ws.Cells.Item(102, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(102, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(102, 16).Value = 1333
$ws.Cells.Item(102, 17).Value = 18
$ws.Cells.Item(102, 18).Value = "Hortaliza"

# New row 103 data (Vega Modelo de Temuco - Pepino dulce)
$ws.Cells.Item(103, 1).Value = 10
$ws.Cells.Item(103, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(103, 3).Value = "La Araucanía"
$ws.Cells.Item(103, 4).Value = 44466
$ws.Cells.Item(103, 5).Value = 9
$ws.Cells.Item(103, 6).Value = 100112043
$ws.Cells.Item(103, 7).Value = "Pepino dulce"
$ws.Cells.Item(103, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(103, 9).Value = "Segunda"
$ws.Cells.Item(103, 10).Value = 30
$ws.Cells.Item(103, 11).Value = 20000
$ws.Cells.Item(103, 12).Value = 20000
$ws.Cells.Item(103, 13).Value = 20000
$ws.Cells.Item(103, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(103, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(103, 16).Value = 1111
$ws.Cells.Item(103, 17).Value = 18
$ws.Cells.Item(103, 18).Value = "Hortaliza"
